# Refresh "paises.xlsx" covid stats sheet with the newer data snapshot
# (countries & provincias Spain update) and bump the "last updated" footer.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Footer timestamp: 12:54 -> 14:11
$ws.Range("A1").Value = 'Datos actualizados a 25 de Junio de 2020 a las 14:11'

# Estados Unidos
$ws.Range("B4").Value = 2463438
$ws.Range("C4").Value = 884
$ws.Range("D4").Value = 1040608
$ws.Range("E4").Value = 1298536
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 13
$ws.Range("H4").Value = 124294

# Brasil
$ws.Range("B5").Value = 1193609
$ws.Range("C5").Value = 1135
$ws.Range("D5").Value = 649908
$ws.Range("E5").Value = 489806
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 21
$ws.Range("H5").Value = 53895

# India
$ws.Range("B7").Value = 474585
$ws.Range("C7").Value = 1600
$ws.Range("D7").Value = 272382
$ws.Range("E7").Value = 187288
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 8
$ws.Range("H7").Value = 14915

# Alemania
$ws.Range("B15").Value = 193281
$ws.Range("C15").Value = 27
$ws.Range("D15").Value = 176800
$ws.Range("E15").Value = 7478

# Catar
$ws.Range("B23").Value = 91838
$ws.Range("C23").Value = 1060
$ws.Range("D23").Value = 74544
$ws.Range("E23").Value = 17188
$ws.Range("F23").Value = 0
$ws.Range("G23").Value = 2
$ws.Range("H23").Value = 106

# Bielorrusia
$ws.Range("B28").Value = 60382
$ws.Range("C28").Value = 437
$ws.Range("D28").Value = 41448
$ws.Range("E28").Value = 18567
$ws.Range("F28").Value = 0
$ws.Range("G28").Value = 5
$ws.Range("H28").Value = 367

# Emiratos Arabes Unidos
$ws.Range("B34").Value = 46563
$ws.Range("C34").Value = 430
$ws.Range("D34").Value = 35165
$ws.Range("E34").Value = 11090
$ws.Range("F34").Value = 0
$ws.Range("G34").Value = 1
$ws.Range("H34").Value = 308

# Row 35/36: Kuwait and Singapur swap order (re-sorted by Casos totales)
$ws.Range("A35").Value = 'Kuwait'
$ws.Range("B35").Value = 42788
$ws.Range("C35").Value = 909
$ws.Range("D35").Value = 33367
$ws.Range("E35").Value = 9082
$ws.Range("F35").Value = 0
$ws.Range("G35").Value = 2
$ws.Range("H35").Value = 339

$ws.Range("A36").Value = 'Singapur'
$ws.Range("B36").Value = 42736
$ws.Range("C36").Value = 113
$ws.Range("D36").Value = 36299
$ws.Range("E36").Value = 6411
$ws.Range("F36").Value = 0
$ws.Range("G36").Value = 0
$ws.Range("H36").Value = 26

$ws.Range("B37").Value = 40415
$ws.Range("C37").Value = 311
$ws.Range("D37").Value = 26382
$ws.Range("E37").Value = 12484
$ws.Range("F37").Value = 0
$ws.Range("G37").Value = 6
$ws.Range("H37").Value = 1549

# Row 101/102: Croacia and Mayotte swap order
$ws.Range("A101").Value = 'Croacia'
$ws.Range("B101").Value = 2483
$ws.Range("C101").Value = 95
$ws.Range("D101").Value = 2149
$ws.Range("E101").Value = 227
$ws.Range("F101").Value = 0
$ws.Range("G101").Value = 0
$ws.Range("H101").Value = 107

$ws.Range("A102").Value = 'Mayotte'
$ws.Range("B102").Value = 2467
$ws.Range("C102").Value = 0
$ws.Range("D102").Value = 2218
$ws.Range("E102").Value = 217
$ws.Range("F102").Value = 0
$ws.Range("G102").Value = 0
$ws.Range("H102").Value = 32

# Rows 111-113: Madagascar moves ahead of Islandia / Lituania
$ws.Range("A111").Value = 'Madagascar'
$ws.Range("B111").Value = 1829
$ws.Range("C111").Value = 42
$ws.Range("D111").Value = 823
$ws.Range("E111").Value = 990
$ws.Range("F111").Value = 0
$ws.Range("G111").Value = 0
$ws.Range("H111").Value = 16

$ws.Range("A112").Value = 'Islandia'
$ws.Range("B112").Value = 1824
$ws.Range("C112").Value = 0
$ws.Range("D112").Value = 1806
$ws.Range("E112").Value = 8
$ws.Range("F112").Value = 0
$ws.Range("G112").Value = 0
$ws.Range("H112").Value = 10

$ws.Range("A113").Value = 'Lituania'
$ws.Range("B113").Value = 1806
$ws.Range("C113").Value = 2
$ws.Range("D113").Value = 1494
$ws.Range("E113").Value = 234
$ws.Range("F113").Value = 0
$ws.Range("G113").Value = 0
$ws.Range("H113").Value = 78

# Hong Kong
$ws.Range("B124").Value = 1194
$ws.Range("C124").Value = 14
$ws.Range("D124").Value = 1088
$ws.Range("E124").Value = 99
$ws.Range("F124").Value = 0
$ws.Range("G124").Value = 1
$ws.Range("H124").Value = 7

# Burkina Faso
$ws.Range("B135").Value = 934
$ws.Range("C135").Value = 15
$ws.Range("D135").Value = 826
$ws.Range("E135").Value = 55

# Row 202/203: Dominica and Fiyi swap order (tied totals)
$ws.Range("A202").Value = 'Dominica'
$ws.Range("A203").Value = 'Fiyi'

# Row 208/209: Islas Malvinas and Groenlandia swap order (tied totals)
$ws.Range("A208").Value = 'Islas Malvinas'
$ws.Range("A209").Value = 'Groenlandia'
